$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "Datos actualizados" timestamp in A1 ---
$ws.Range("A1").Value = "Datos actualizados a 18 de Junio de 2020 a las 21:46"

# --- Country name swaps (adjacent rows whose rank order changed) ---
$ws.Range("A23").Value = "Sudafrica"
$ws.Range("A24").Value = "China"
$ws.Range("A50").Value = "Israel"
$ws.Range("A51").Value = "Barein"
$ws.Range("A105").Value = "Costa Rica"
$ws.Range("A106").Value = "Mali"
$ws.Range("A154").Value = "Zimbabue"
$ws.Range("A155").Value = "Taiwan"
$ws.Range("A208").Value = "Islas Turcas y Caicos"
$ws.Range("A209").Value = "Santa Sede"
$ws.Range("A213").Value = "Islas Virgenes Britanicas"
$ws.Range("A214").Value = "Papua Nueva Guinea"

# --- Updated statistic values ---
$ws.Range("B7").Value = 381091
$ws.Range("C7").Value = 13827
$ws.Range("D7").Value = 205182
$ws.Range("E7").Value = 163305
$ws.Range("G7").Value = 342
$ws.Range("H7").Value = 12604
$ws.Range("B21").Value = 100148
$ws.Range("C21").Value = 295
$ws.Range("E21").Value = 29404
$ws.Range("B23").Value = 83890
$ws.Range("C23").Value = 3478
$ws.Range("D23").Value = 44920
$ws.Range("E23").Value = 37233
$ws.Range("G23").Value = 63
$ws.Range("H23").Value = 1737
$ws.Range("B24").Value = 83293
$ws.Range("C24").Value = 28
$ws.Range("D24").Value = 78394
$ws.Range("E24").Value = 265
$ws.Range("H24").Value = 4634
$ws.Range("B50").Value = 19998
$ws.Range("C50").Value = 215
$ws.Range("D50").Value = 15518
$ws.Range("E50").Value = 4177
$ws.Range("G50").Value = 0
$ws.Range("H50").Value = 303
$ws.Range("B51").Value = 19961
$ws.Range("C51").Value = 0
$ws.Range("D51").Value = 14185
$ws.Range("E51").Value = 5723
$ws.Range("G51").Value = 4
$ws.Range("H51").Value = 53
$ws.Range("B102").Value = 2137
$ws.Range("C102").Value = 17
$ws.Range("E102").Value = 452
$ws.Range("B104").Value = 1946
$ws.Range("C104").Value = 22
$ws.Range("E104").Value = 514
$ws.Range("B105").Value = 1939
$ws.Range("C105").Value = 68
$ws.Range("D105").Value = 937
$ws.Range("E105").Value = 990
$ws.Range("H105").Value = 12
$ws.Range("B106").Value = 1906
$ws.Range("C106").Value = 16
$ws.Range("D106").Value = 1192
$ws.Range("E106").Value = 607
$ws.Range("H106").Value = 107
$ws.Range("B154").Value = 463
$ws.Range("C154").Value = 62
$ws.Range("D154").Value = 63
$ws.Range("E154").Value = 396
$ws.Range("H154").Value = 4
$ws.Range("B155").Value = 446
$ws.Range("C155").Value = 1
$ws.Range("D155").Value = 434
$ws.Range("E155").Value = 5
$ws.Range("H155").Value = 7
$ws.Range("B171").Value = 166
$ws.Range("C171").Value = 11
$ws.Range("E171").Value = 94
$ws.Range("G171").Value = 1
$ws.Range("H171").Value = 8
$ws.Range("D208").Value = 11
$ws.Range("H208").Value = 1
$ws.Range("D209").Value = 12
$ws.Range("H209").Value = 0
$ws.Range("D213").Value = 7
$ws.Range("H213").Value = 1
$ws.Range("D214").Value = 8
$ws.Range("H214").Value = 0
